$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("GRUPO_SINPAR")
$ws2 = $wb.Worksheets.Item("GRUPO_MAYORISTAS")

# Shift existing data (rows 2..9) down to rows 3..10 manually, bottom-up,
# to avoid Excel's Insert() copying formatting down from the header row.
for ($r = 9; $r -ge 2; $r--) {
    $srcVal = $ws2.Cells.Item($r, 1).Value()
    $ws2.Cells.Item($r + 1, 1).Value = $srcVal
}

# Set the new cell's value (same shared string as GRUPO_SINPAR!A2 -> "evol0088")
$ws2.Range("A2").Value = "evol0088"

# Apply left-horizontal alignment style to the new cell
$ws2.Range("A2").HorizontalAlignment = -4131  # xlLeft

# Update selections to match target state
$ws1.Range("A2").Select()
$ws2.Range("C5").Select()
